# Add 2022-Q4 data
# -----------------
# 1) Duplicate the existing "2022-Q3" sheet (preserves all cell styles),
#    placing the copy immediately before it, then rename the copy to
#    "2022-Q4" and update its figures to the new quarter's numbers.
# 2) Insert the new quarter into the "总计" (totals) summary sheet as the
#    new first data row, shifting the older rows down and renumbering the
#    index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write $text into $range as a genuine text value (not a number),
# without touching the cell's style/number-format - mirrors how the
# existing sheet stores figures such as "18.69" / "0.10" as text.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($ws, $cellRef, $text)
    $helper = $ws.Range("Z100")
    $helper.Formula = "=""" + $text + """"
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $helper.ClearContents()
}

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q4" sheet from a copy of "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Update the fund-scale / position figures for the new quarter
Set-TextValue $q4 "D2" "17.24"
Set-TextValue $q4 "E2" "93.84"
Set-TextValue $q4 "F2" "3.09"
Set-TextValue $q4 "G2" "0.5327"
$q4.Range("H2").Value = 9

Set-TextValue $q4 "D3" "1.56"
Set-TextValue $q4 "E3" "93.84"
Set-TextValue $q4 "F3" "3.09"
Set-TextValue $q4 "G3" "0.0482"
$q4.Range("H3").Value = 9

Set-TextValue $q4 "D4" "0.20"
Set-TextValue $q4 "E4" "93.84"
Set-TextValue $q4 "F4" "3.09"
Set-TextValue $q4 "G4" "0.0062"
$q4.Range("H4").Value = 9

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend column A's formatting down to the new last row (6)
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

# Shift the existing four data rows down by one (bottom-up so each source
# cell is read before it gets overwritten), renumbering the index column.
$total.Range("A6").Value = 4
$total.Range("B6").Value = $total.Range("B5").Value2
$total.Range("C6").Value = $total.Range("C5").Value2
$total.Range("D6").Value = $total.Range("D5").Value2

$total.Range("A5").Value = 3
$total.Range("B5").Value = $total.Range("B4").Value2
$total.Range("C5").Value = $total.Range("C4").Value2
$total.Range("D5").Value = $total.Range("D4").Value2

$total.Range("A4").Value = 2
$total.Range("B4").Value = $total.Range("B3").Value2
$total.Range("C4").Value = $total.Range("C3").Value2
$total.Range("D4").Value = $total.Range("D3").Value2

$total.Range("A3").Value = 1
$total.Range("B3").Value = $total.Range("B2").Value2
$total.Range("C3").Value = $total.Range("C2").Value2
$total.Range("D3").Value = $total.Range("D2").Value2

# New first data row: 2022-Q4
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.59
